# Updated cryptos list on Wed May 29 10:59:02 UTC 2024 with GitHub Actions
#
# This script mirrors a refresh of the scraped coinranking.com price table:
# most rows get updated Price (col D) and/or Volume(1h) (col E) values,
# and rows 48/49 swap which coin (OKB / FLOKI) occupies which row.
#
# Price cells are plain scraped text (e.g. "67.893.45", "1.00"), so when the
# new value would otherwise be auto-recognised by Excel as a number we force
# the cell's number format to Text ("@") first, keeping the stored cell type
# as text/string, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    # Only force text formatting when the literal would otherwise be
    # re-interpreted by Excel as a number (single-dot decimals etc.) -
    # values with two dots ("67.952.31"), letters, or other
    # non-numeric-looking text are already safe to assign directly and this
    # keeps us from touching their style unnecessarily.
    if ($Text -match '^[0-9]+(\.[0-9]+)?$') {
        $Range.NumberFormat = "@"
    }
    $Range.Value = $Text
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "67.952.31"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.820.90"
$ws.Range("E3").Value = "  -2.13%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "600.48"
$ws.Range("E5").Value = "  -0.47%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "169.91"
$ws.Range("E6").Value = "  +0.14%  "

# Row 7 - LidoStakedEther
Set-TextValue $ws.Range("D7") "3.821.33"
$ws.Range("E7").Value = "  -2.14%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.166"
$ws.Range("E10").Value = "  -0.87%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "6.52"
$ws.Range("E11").Value = "  +0.85%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.463"
$ws.Range("E12").Value = "  +0.61%  "

# Row 13 - ShibaInu
Set-TextValue $ws.Range("D13") "0.0000278"
$ws.Range("E13").Value = "  +9.00%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "37.14"
$ws.Range("E14").Value = "  -0.14%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "4.463.45"

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "3.801.75"
$ws.Range("E16").Value = "  -2.38%  "

# Row 17 - Chainlink
$ws.Range("E17").Value = "  +3.06%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "68.003.67"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "7.46"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.19%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "10.90"
$ws.Range("E21").Value = "  +0.33%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "469.97"
$ws.Range("E22").Value = "  -0.79%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +0.25%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  -9.14%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "83.66"

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +2.51%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +0.21%  "

# Row 28 - RenderToken
Set-TextValue $ws.Range("D28") "10.39"
$ws.Range("E28").Value = "  +3.74%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.09%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.67%  "

# Row 31 - WrappedeETH
Set-TextValue $ws.Range("D31") "3.967.62"
$ws.Range("E31").Value = "  -2.08%  "

# Row 32 - NEARProtocol
Set-TextValue $ws.Range("D32") "7.77"
$ws.Range("E32").Value = "  -1.48%  "

# Row 33 - ImmutableX
Set-TextValue $ws.Range("D33") "2.29"
$ws.Range("E33").Value = "  -1.40%  "

# Row 34 - EthereumClassic
Set-TextValue $ws.Range("D34") "30.81"
$ws.Range("E34").Value = "  -2.54%  "

# Row 35 - Aptos
Set-TextValue $ws.Range("D35") "9.39"
$ws.Range("E35").Value = "  -0.60%  "

# Row 36 - RenzoRestakedETH
Set-TextValue $ws.Range("D36") "3.786.03"
$ws.Range("E36").Value = "  -2.28%  "

# Row 37 - dogwifhat
Set-TextValue $ws.Range("D37") "3.90"
$ws.Range("E37").Value = "  +4.49%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  +1.53%  "

# Row 39 - Filecoin
Set-TextValue $ws.Range("D39") "6.01"
$ws.Range("E39").Value = "  +1.31%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -1.35%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  -2.49%  "

# Row 42 - FirstDigitalUSD
Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  +1.89%  "

# Row 44 - Cosmos
$ws.Range("E44").Value = "  +2.29%  "

# Row 45 - USDe: unchanged

# Row 46 - Stacks
Set-TextValue $ws.Range("D46") "1.98"
$ws.Range("E46").Value = "  -1.39%  "

# Row 47 - Bittensor
Set-TextValue $ws.Range("D47") "412.70"
$ws.Range("E47").Value = "  -3.97%  "

# Rows 48/49 swap: OKB and FLOKI trade places (with refreshed price/volume)
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D48") "0.000288"
$ws.Range("E48").Value = "  -5.19%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D49") "46.58"
$ws.Range("E49").Value = "  -1.47%  "

# Row 50 - Monero
Set-TextValue $ws.Range("D50") "142.78"
$ws.Range("E50").Value = "  -0.97%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -0.17%  "
